# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Zapallo italiano" (Agrícola del
# Norte S.A. de Arica) at row 372, pushing the existing rows 372:394 down to
# 373:395.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 372:394 down to 373:395, leaving a fresh blank row 372.
$ws.Rows(372).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(372, 1).Value = 1
$ws.Cells.Item(372, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(372, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(372, 4).Value = 44783
$ws.Cells.Item(372, 5).Value = 15
$ws.Cells.Item(372, 6).Value = 100112032
$ws.Cells.Item(372, 7).Value = "Zapallo italiano"
$ws.Cells.Item(372, 8).Value = "Huracán"
$ws.Cells.Item(372, 9).Value = "Segunda"
$ws.Cells.Item(372, 10).Value = 100
$ws.Cells.Item(372, 11).Value = 14000
$ws.Cells.Item(372, 12).Value = 15000
$ws.Cells.Item(372, 13).Value = 14500
$ws.Cells.Item(372, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(372, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(372, 16).Value = 145
$ws.Cells.Item(372, 17).Value = 100
$ws.Cells.Item(372, 18).Value = "Hortaliza"

# Preserve the date-formatted number format already used by the other rows'
# date column (column D) for the newly inserted cell.
$ws.Range("D372").NumberFormat = "YYYY-MM-DD HH:MM:SS"
